$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix the shared string text for the "Restaurants, cafe, etc." label (cell Z1)
$ws.Range("Z1").Value = "Restaurant, Cafes etc."

# 2. Replace the #NUM! errors in column F (rows 2-64) with numeric 0
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# 3. Append a new data row (row 65) with the supplied values
$row65 = @{
    "A65" = 0.00000557903116981811
    "B65" = 0.000045912798512734
    "C65" = 0.00087540118877356
    "D65" = 0.00000577884284510534
    "E65" = 0.0000112356501341112
    "F65" = 0
    "G65" = 0.0000038506029333242
    "H65" = 0.000000518185406244657
    "I65" = 0.0000376224427349759
    "J65" = 0.000000282184265118779
    "K65" = 0.000083244054027958
    "L65" = 0.0000412425375561132
    "M65" = 0
    "N65" = 0.0000642639559525381
    "O65" = 0.0000156912605590989
    "P65" = 0.000123904701849719
    "Q65" = 0.000000964672142964161
    "R65" = 0.0000182307615169893
    "S65" = 0.000514905681648989
    "T65" = 0.000154976205688869
    "U65" = 0.00000789375636598898
    "V65" = 0.000000532990703565931
    "W65" = 0.0000000409598026189442
    "X65" = 0.000211222118858404
    "Y65" = 0.0000957596663659824
    "Z65" = 0.00157260308073657
    "AA65" = 0.0000797302005192069
    "AB65" = 0.000138932552975629
    "AC65" = 0.000192787662946777
    "AD65" = 0.0000455766242789465
    "AE65" = 0
    "AF65" = 0.00000117583053310462
    "AG65" = -0.000211865313747471
    "AH65" = 0.00000507394050056123
    "AI65" = 0
    "AJ65" = 0.000001181667060156
    "AK65" = 0.0000318156054330826
}

foreach ($addr in $row65.Keys) {
    $ws.Range($addr).Value = $row65[$addr]
}
